$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("Q15").Value = 851.76
$ws1.Range("Q19").Value = "1 de 17"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F15").Value = 2100.7
$ws2.Range("F19").Value = 32556.53

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D14").Value = 851.76
$ws3.Range("E14").Value = -368.76
$ws3.Range("F14").Value = 1.763478260869565

$ws3.Range("D19").Value = 32556.53
$ws3.Range("E19").Value = 14662.77386304603
$ws3.Range("F19").Value = 0.6894750099329362
